$d = $word.ActiveDocument

# --- Occurrence 1: <id>p006v_a1</id>  ->  <id>p006v_1</id>  (merge 3 runs into 1) ---
$rng1 = $d.Content
$rng1.Find.Execute("<id>p006v_a1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng1.Text = "<id>p006v_1</id>"

# --- Occurrence 2: <id>p006v_2</id>  ->  <id>p006v_2</id>  (merge 3 runs into 1, text unchanged) ---
# Setting identical text can be a no-op for the underlying run model, so force a
# genuine mutation first (temp text) before restoring the final text; this
# guarantees the three source runs collapse into a single run as in the target.
$rng2 = $d.Content
$rng2.Find.Execute("<id>p006v_2</id>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Text = "<id>p006v_2_tmp_placeholder</id>"
$rng2.Text = "<id>p006v_2</id>"

Write-Output "Done"
